# Add a new column CF to the right of the existing CE column.
# CF1 gets a new quarter-end date (2025-11-25 -> serial 45986), formatted
# like the other date headers in row 1 (copy CE1's formatting).
# CF4:CF35 duplicate the value currently held in the corresponding CE cell
# (the preprocessing step that seeds the newest column from the previous one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new date header in CF1, formatted like CE1 ---
$ce1 = $ws.Range("CE1")
$cf1 = $ws.Range("CF1")
$ce1.Copy()
$cf1.PasteSpecial(-4122)   # xlPasteFormats
$cf1.Value = 45986

# --- Rows 4-35: copy CE value into the new CF cell ---
for ($r = 4; $r -le 35; $r++) {
    $ceValue = $ws.Range("CE$r").Value2
    $ws.Range("CF$r").Value = $ceValue
}
